# "Verificação do que Está Feito"
# Adds an "Adicionar Vencedor" test row, splits the "Associar Jogador a uma
# Equipa" row into two (club + team-of-club), adds an "Associar Clube
# Vencedor" row, appends a new "//UTILIZADORES" test section, and marks
# each existing test row as OK (column B) or Not Okay (column C).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1) Insert "Adicionar Vencedor" row right after
#    "Associar a um Formato de Competição" (before "//ASSOCIAÇÕES DE FUTEBOL")
# ---------------------------------------------------------------------
$ws.Rows.Item(7).Insert()
$ws.Range("A7").Value = "Adicionar Vencedor"

# ---------------------------------------------------------------------
# 2) Split "Associar Jogador a uma Equipa" (old row 31, now row 32 after
#    the insert above) into two rows: "...a um Clube" and
#    "...a uma Equipa do respetivo Clube"
# ---------------------------------------------------------------------
$ws.Rows.Item(32).Insert()
$ws.Range("A32").Value = "Associar Jogador a um Clube"
$ws.Range("A33").Value = "Associar Jogador a uma Equipa do respetivo Clube"
$ws.Rows.Item(33).RowHeight = 28.8

# ---------------------------------------------------------------------
# 3) Insert "Associar Clube Vencedor" row right after
#    "Associar Jogo a um Estádio" (before "//ESTATÍSTICAS JOGO")
# ---------------------------------------------------------------------
$ws.Rows.Item(42).Insert()
$ws.Range("A42").Value = "Associar Clube Vencedor"

# ---------------------------------------------------------------------
# 4) Append the new "//UTILIZADORES" test section at the bottom
# ---------------------------------------------------------------------
$ws.Range("A52").Value = "//UTILIZADORES"
$ws.Range("A53").Value = "Criar Conta"
$ws.Range("A54").Value = "Conta Admin Criada tem de ser Ativada"
$ws.Range("A55").Value = "Iniciar Sessão"
$ws.Range("A56").Value = "Terminar Sessão"
$ws.Range("A57").Value = "Utilizador tem Clubes Favoritos"

# ---------------------------------------------------------------------
# 5) Mark header/section rows ("//...") with "//" in both B and C, and
#    every regular test row with "X" in either B (OK) or C (Not Okay).
# ---------------------------------------------------------------------
$okRows       = @(3,4,5,6,9,10,11,13,14,15,16,17,19,20,21,22,24,25,26,28,29,30,31,32,33,35,37,38,39,40,41,44,45,46,47,48,49,50,51,53,54,55,56,57)
$notOkayRows  = @(7,36,42)
$headerRows   = @(2,8,12,18,23,27,34,43,52)

foreach ($r in $headerRows) {
    $ws.Range("B$r").Value = "//"
    $ws.Range("C$r").Value = "//"
}
foreach ($r in $okRows) {
    $ws.Range("B$r").Value = "X"
}
foreach ($r in $notOkayRows) {
    $ws.Range("C$r").Value = "X"
}

# ---------------------------------------------------------------------
# 6) Restore the selection / scroll state recorded in the workbook
# ---------------------------------------------------------------------
$ws.Range("B47").Select()
